$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# Summary sheet
# -----------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Activate()

$wsSummary.Range("B2").Value = 848.21
$wsSummary.Range("E2").NumberFormat = "#,##0.00"
$wsSummary.Range("E2").Value = 9151.79
$wsSummary.Range("F2").Value = 857.01

$wsSummary.Range("A3").Value = 555.57
$wsSummary.Range("E3").Value = 503.78
$wsSummary.Range("F3").Value = 42.99

$wsSummary.Range("A7:XFD15").Select() | Out-Null

# -----------------------------------------------------------------------
# Repayment schedule sheet
# -----------------------------------------------------------------------
$wsRepay = $wb.Worksheets.Item("Repayment schedule")
$wsRepay.Activate()

# Bring in new O3:O15 / P2 cells with the same formatting as their
# already-present neighbours (numFmtId 0, vertical-center + wrap) before
# writing values into them.
$wsRepay.Range("O2").Copy()
$wsRepay.Range("O3:O15").PasteSpecial(-4122)
$wsRepay.Range("P2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 4 is a disbursement marker row (like row 2) - O4 stays empty, only
# gets the inherited formatting from the PasteSpecial above. Every other
# row in 3..15 gets an explicit 0.
foreach ($r in @(3,5,6,7,8,9,10,11,12,13,14,15)) {
    $wsRepay.Range("O$r").Value = 0
}

# Row 3
$wsRepay.Range("F3").Value = 848.21
$wsRepay.Range("G3").NumberFormat = "#,##0.00"
$wsRepay.Range("G3").Value = 4151.79
$wsRepay.Range("K3").Value = 900
$wsRepay.Range("L3").Value = 900

# Row 5
$wsRepay.Range("F5").Value = 857.01
$wsRepay.Range("G5").NumberFormat = "#,##0.00"
$wsRepay.Range("G5").Value = 8294.78
$wsRepay.Range("H5").Value = 42.99
$wsRepay.Range("K5").Value = 900
$wsRepay.Range("P5").Value = 900

# Row 6
$wsRepay.Range("F6").Value = 808.67
$wsRepay.Range("G6").NumberFormat = "#,##0.00"
$wsRepay.Range("G6").Value = 7486.11
$wsRepay.Range("H6").Value = 91.33
$wsRepay.Range("K6").Value = 900
$wsRepay.Range("P6").Value = 900

# Row 7
$wsRepay.Range("F7").Value = 825.14
$wsRepay.Range("G7").NumberFormat = "#,##0.00"
$wsRepay.Range("G7").Value = 6660.97
$wsRepay.Range("H7").Value = 74.86
$wsRepay.Range("K7").Value = 900
$wsRepay.Range("P7").Value = 900

# Row 8
$wsRepay.Range("F8").Value = 833.39
$wsRepay.Range("G8").NumberFormat = "#,##0.00"
$wsRepay.Range("G8").Value = 5827.58
$wsRepay.Range("H8").Value = 66.61
$wsRepay.Range("K8").Value = 900
$wsRepay.Range("P8").Value = 900

# Row 9
$wsRepay.Range("F9").Value = 841.72
$wsRepay.Range("G9").NumberFormat = "#,##0.00"
$wsRepay.Range("G9").Value = 4985.86
$wsRepay.Range("H9").Value = 58.28
$wsRepay.Range("K9").Value = 900
$wsRepay.Range("P9").Value = 900

# Row 10
$wsRepay.Range("F10").Value = 850.14
$wsRepay.Range("G10").NumberFormat = "#,##0.00"
$wsRepay.Range("G10").Value = 4135.72
$wsRepay.Range("H10").Value = 49.86
$wsRepay.Range("K10").Value = 900
$wsRepay.Range("P10").Value = 900

# Row 11
$wsRepay.Range("F11").Value = 858.64
$wsRepay.Range("G11").NumberFormat = "#,##0.00"
$wsRepay.Range("G11").Value = 3277.08
$wsRepay.Range("H11").Value = 41.36
$wsRepay.Range("K11").Value = 900
$wsRepay.Range("P11").Value = 900

# Row 12
$wsRepay.Range("F12").Value = 867.23
$wsRepay.Range("G12").NumberFormat = "#,##0.00"
$wsRepay.Range("G12").Value = 2409.85
$wsRepay.Range("H12").Value = 32.77
$wsRepay.Range("K12").Value = 900
$wsRepay.Range("P12").Value = 900

# Row 13
$wsRepay.Range("F13").Value = 875.9
$wsRepay.Range("G13").NumberFormat = "#,##0.00"
$wsRepay.Range("G13").Value = 1533.95
$wsRepay.Range("H13").Value = 24.1
$wsRepay.Range("K13").Value = 900
$wsRepay.Range("P13").Value = 900

# Row 14 (G14 keeps its original style, only the value changes)
$wsRepay.Range("F14").Value = 884.66
$wsRepay.Range("G14").Value = 649.29
$wsRepay.Range("H14").Value = 15.34
$wsRepay.Range("K14").Value = 900
$wsRepay.Range("P14").Value = 900

# Row 15
$wsRepay.Range("F15").Value = 649.29
$wsRepay.Range("H15").Value = 6.28
$wsRepay.Range("K15").Value = 655.57
$wsRepay.Range("P15").Value = 655.57

$wsRepay.Range("A16:XFD16").Select() | Out-Null

# -----------------------------------------------------------------------
# Transactions sheet
# -----------------------------------------------------------------------
$wsTrans = $wb.Worksheets.Item("Transactions")
$wsTrans.Activate()

$wsTrans.Range("A2").Value = 90
$wsTrans.Range("J2").NumberFormat = "#,##0.00"
$wsTrans.Range("J2").Value = 9151.79

$wsTrans.Range("A3").Value = 89
$wsTrans.Range("E3").Value = 900
$wsTrans.Range("F3").Value = 848.21
$wsTrans.Range("J3").NumberFormat = "#,##0.00"
$wsTrans.Range("J3").Value = 4151.79

$wsTrans.Range("A4").Value = 88

$wsTrans.Range("D3").Select() | Out-Null
